$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9974.99936182774
$ws.Range("F2").Value = -5.53154655365855

$ws.Range("C3").Value = 9345.74916030489
$ws.Range("F3").Value = 198.709128153699

$ws.Range("C4").Value = 6691.8283293934
$ws.Range("F4").Value = 76.1191273590696

$ws.Range("C5").Value = 6566.22568372998
$ws.Range("F5").Value = 68.2775932200327

$ws.Range("C6").Value = 6476.22354814191
$ws.Range("F6").Value = 65.9417813673001

$ws.Range("C7").Value = 9651.47316146217
$ws.Range("F7").Value = 257.858132013451

$ws.Range("C9").Value = 9575.14607243933
$ws.Range("F9").Value = 254.677836637499

$ws.Range("C10").Value = 8773.6790228366
$ws.Range("F10").Value = 221.283376237386

$ws.Range("C11").Value = 8662.23805850486
$ws.Range("F11").Value = 199.857012994936

$ws.Range("C12").Value = 8489.52525201526
$ws.Range("F12").Value = 192.280709421214

$ws.Range("C13").Value = 9090.68588180193
$ws.Range("F13").Value = 245.85832935733

$ws.Range("C14").Value = 9313.41589184908
$ws.Range("F14").Value = 255.138746442628

$ws.Range("C15").Value = 9409.05495233866
$ws.Range("F15").Value = 259.12370729636
